# Trade #40 closed at 2026-02-17 08:33:21 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet - update aggregate metrics
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.11   # Total P&L %
$summary.Range("B6").Value = 40      # Total Trades
$summary.Range("B9").Value = 35      # Win Rate %

# ---------------------------------------------------------------
# 2) Strategy Status sheet - update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 40       # Trades
$status.Range("G4").Value = 35       # Win Rate %

# ---------------------------------------------------------------
# Helper: write one new trade row (values only, columns A..Q) onto a
# given worksheet at the given row number. Column B holds a
# "yyyy-mm-dd"-looking string which Excel's smart-typing would
# otherwise coerce into a real date serial; force it to stay plain
# text (matching the rest of the column) and then restore the
# "Normal" style so no stray per-cell formatting is left behind.
# ---------------------------------------------------------------
function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 40
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value = "08:33:15"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.43
    $ws.Cells.Item($row, 7).Value = 0.43
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 99.77
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

# ---------------------------------------------------------------
# 3) All Trades sheet - append new trade row 41
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 41

# ---------------------------------------------------------------
# 4) MarketMaking sheet - append the same new trade row 41
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $mm 41
